$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '69.163.94'
$ws.Range("E2").Value = '  -0.35%  '
$ws.Range("D3").Value = '3.669.84'
$ws.Range("E3").Value = '  -0.42%  '
$ws.Range("E4").Value = '  +0.06%  '
$ws.Range("D5").Value = '674.40'
$ws.Range("E5").Value = '  -1.11%  '
$ws.Range("D6").Value = '157.23'
$ws.Range("E6").Value = '  -3.23%  '
$ws.Range("E7").Value = '  +0.05%  '
$ws.Range("D8").Value = '0.492'
$ws.Range("E8").Value = '  -1.54%  '
$ws.Range("D9").Value = '0.145'
$ws.Range("E9").Value = '  -1.82%  '
$ws.Range("D10").Value = '6.92'
$ws.Range("E10").Value = '  -5.75%  '
$ws.Range("D11").Value = '0.435'
$ws.Range("E11").Value = '  -2.50%  '
$ws.Range("D12").Value = '0.0000230'
$ws.Range("E12").Value = '  -3.69%  '
$ws.Range("D13").Value = '4.288.74'
$ws.Range("E13").Value = '  -0.46%  '
$ws.Range("D14").Value = '32.14'
$ws.Range("E14").Value = '  -3.97%  '
$ws.Range("D15").Value = '3.672.68'
$ws.Range("E15").Value = '  -0.48%  '
$ws.Range("D16").Value = '69.150.04'
$ws.Range("E16").Value = '  -0.37%  '
$ws.Range("E17").Value = '  +0.92%  '
$ws.Range("D18").Value = '15.99'
$ws.Range("E18").Value = '  -1.59%  '
$ws.Range("D19").Value = '6.39'
$ws.Range("E19").Value = '  -3.47%  '
$ws.Range("D20").Value = '466.50'
$ws.Range("E20").Value = '  -3.20%  '
$ws.Range("D21").Value = '9.93'
$ws.Range("E21").Value = '  +0.45%  '
$ws.Range("D22").Value = '0.646'
$ws.Range("E22").Value = '  -2.93%  '
$ws.Range("D23").Value = '79.83'
$ws.Range("E23").Value = '  -0.58%  '
$ws.Range("D24").Value = '3.818.63'
$ws.Range("E24").Value = '  -0.34%  '
$ws.Range("E25").Value = '  +0.00%  '
$ws.Range("B26").Value = 'InternetComputer(DFINITY)'
$ws.Range("C26").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D26").Value = '10.85'
$ws.Range("E26").Value = '  -5.70%  '
$ws.Range("B27").Value = 'PEPE'
$ws.Range("C27").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D27").Value = '0.0000119'
$ws.Range("E27").Value = '  -8.43%  '
$ws.Range("D28").Value = '8.95'
$ws.Range("E28").Value = '  -6.71%  '
$ws.Range("D29").Value = '2.66'
$ws.Range("E29").Value = '  -2.42%  '
$ws.Range("D30").Value = '1.72'
$ws.Range("E30").Value = '  -6.61%  '
$ws.Range("E31").Value = '  -3.84%  '
$ws.Range("E32").Value = '  +0.00%  '
$ws.Range("D33").Value = '26.82'
$ws.Range("E33").Value = '  -1.02%  '
$ws.Range("E34").Value = '  -5.50%  '
$ws.Range("D35").Value = '3.660.50'
$ws.Range("E35").Value = '  +0.20%  '
$ws.Range("D36").Value = '0.160'
$ws.Range("E36").Value = '  -4.39%  '
$ws.Range("D37").Value = '8.10'
$ws.Range("E37").Value = '  -4.72%  '
$ws.Range("E38").Value = '  -3.31%  '
$ws.Range("E39").Value = '  +0.00%  '
$ws.Range("E40").Value = '  -0.01%  '
$ws.Range("D41").Value = '2.20'
$ws.Range("E41").Value = '  -2.51%  '
$ws.Range("D42").Value = '0.0894'
$ws.Range("E42").Value = '  -4.30%  '
$ws.Range("D43").Value = '172.02'
$ws.Range("E43").Value = '  +7.21%  '
$ws.Range("D44").Value = '0.939'
$ws.Range("E44").Value = '  -1.73%  '
$ws.Range("D45").Value = '47.47'
$ws.Range("E45").Value = '  -1.85%  '
$ws.Range("D46").Value = '0.000274'
$ws.Range("E46").Value = '  -5.27%  '
$ws.Range("D47").Value = '2.65'
$ws.Range("E47").Value = '  -7.08%  '
$ws.Range("B48").Value = 'ONDO'
$ws.Range("C48").Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range("D48").Value = '1.27'
$ws.Range("E48").Value = '  -6.50%  '
$ws.Range("B49").Value = 'InjectiveProtocol'
$ws.Range("C49").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D49").Value = '27.09'
$ws.Range("E49").Value = '  -9.71%  '
$ws.Range("D50").Value = '1.07'
$ws.Range("E50").Value = '  -2.65%  '
$ws.Range("D51").Value = '7.74'
$ws.Range("E51").Value = '  -3.58%  '
